$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(136).Insert()
$ws.Range("A136").Value = 7
$ws.Range("B136").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C136").Value = "Ñuble"
$ws.Range("D136").Value = 44603
$ws.Range("E136").Value = 16
$ws.Range("F136").Value = 100112023
$ws.Range("G136").Value = "Brócoli"
$ws.Range("H136").Value = "Sin especificar"
$ws.Range("I136").Value = "Primera"
$ws.Range("J136").Value = 160
$ws.Range("K136").Value = 700
$ws.Range("L136").Value = 750
$ws.Range("M136").Value = 725
$ws.Range("N136").Value = "`$/unidad"
$ws.Range("O136").Value = "Provincia de Diguillín"
$ws.Range("P136").Value = 725
$ws.Range("Q136").Value = 1
$ws.Range("R136").Value = "Hortaliza"
